$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.773.30"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "2.247.17"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'303.90"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'94.91"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.486"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'34.56"
$ws.Range("E10").Value = "  +4.47%  "
$ws.Range("D11").Value = "'0.0786"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "'6.74"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "2.597.69"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "2.255.23"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "'0.785"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "41.660.16"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "'12.23"
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").Value = "'5.91"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "'68.00"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "'235.88"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'1.90"
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("D27").Value = "'23.45"
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("D28").Value = "'35.99"
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("D29").Value = "'2.11"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").Value = "'9.37"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").Value = "'159.50"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'5.15"
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("D35").Value = "'0.0729"
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.38"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'16.81"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "'0.103"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "'3.95"
$ws.Range("E41").Value = "  -3.38%  "
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "1.957.03"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("D46").Value = "'9.84"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("D47").Value = "'2.88"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").Value = "'52.43"
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("D49").Value = "'71.62"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").Value = "'90.54"
$ws.Range("E51").Value = "  -1.22%  "
